# "Updates to MN files - IO data, grid battery storage, MCGLT, other files
#  w/ no value change - most recent Jun's export"
#
# For this workbook (Share of Capital and OM Spending by ISIC Code /
# SoCaOMSbRIC) the substantive edits are:
#
#   1. "About" sheet notes text is reworded/condensed:
#        - the note that used to read "...industry energy efficiency,"
#          (continuing onto a separate "CCS, and process emissions
#          policies..." row) becomes a single, self-contained sentence
#          ending "...industry energy efficiency and CCS."
#        - the old "CCS, and process emissions..." row is removed
#        - the old "but in many instances..."/"assume the same spending
#          categorization." rows collapse into a single new note about
#          Process Emissions policies being governed elsewhere
#        - the trailing "assume the same spending categorization." row
#          is removed
#
#   2. "SoCaOMSbRIC" sheet: the combined "ISIC 20T21" column header is
#      split into two separate columns, "ISIC 20" and "ISIC 21" (with a
#      0 data value each), shifting every later ISIC column one slot to
#      the right.
#
# (The shared-strings re-numbering, "ISIC 20T21" removal, and the
# re-indexing visible on the CCS Data sheet are all side effects of
# these two edits and fall out automatically once they're applied.)

$wb = $excel.ActiveWorkbook

# --- 1. "About" sheet: update / collapse the notes text -------------------
$about = $wb.Worksheets.Item("About")

$about.Range("A21").Value2 = "to cover the allocation of capital and OM expensies for industry energy efficiency and CCS."
$about.Range("A22").ClearContents()
$about.Range("A23").Value2 = "A different input variable governs the breakdown of revenues due to Process Emissions policies."
$about.Range("A24").ClearContents()

# --- 2. "SoCaOMSbRIC" sheet: split the "ISIC 20T21" column -----------------
$main = $wb.Worksheets.Item("SoCaOMSbRIC")

# Column K currently holds "ISIC 20T21"; insert a new column after it so
# the combined header can become two separate ones ("ISIC 20"/"ISIC 21"),
# pushing ISIC 22, 23, ... one column to the right.
$main.Columns("L:L").Insert()

$main.Range("K1").Value2 = "ISIC 20"
$main.Range("L1").Value2 = "ISIC 21"
$main.Range("K2").Value2 = 0
$main.Range("L2").Value2 = 0
